# ---------------------------------------------------------------------------
# dataset_tracking.xlsx edit
#
# Adds a second dataset (DS002) to every tracking sheet, alongside small
# corrections to the existing DS001 row (it was mis-tagged as "continuous"
# features / 20C-0D when it should have been "discrete" / 0C-20D, and the
# config filename + creation date + checksum + size were wrong too).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Dataset Registry"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Dataset Registry")

$ws.Cells.Item(2, 2).Value = "dataset-config-001.yml"
$ws.Cells.Item(2, 3).Value = "n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
# Keep this as literal text (not an auto-converted date serial), matching
# the original cell which stored the date as a plain string.
$ws.Cells.Item(2, 4).Value = "'2025-07-02"

$ws.Cells.Item(3, 1).Value = "DS002"
$ws.Cells.Item(3, 2).Value = "dataset-config-002.yml"
$ws.Cells.Item(3, 3).Value = "n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
$ws.Cells.Item(3, 4).Value = "'2025-07-02"
$ws.Cells.Item(3, 5).Value = "Complete"
$ws.Cells.Item(3, 6).Value = "1,000,000 samples, 20 features, linear target"

# ---------------------------------------------------------------------------
# Sheet 2: "Configuration Details"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Configuration Details")

$ws.Cells.Item(2, 2).Value = "dataset-config-001.yml"
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(2, 7).Value = 20

$ws.Cells.Item(3, 1).Value = "DS002"
$ws.Cells.Item(3, 2).Value = "dataset-config-002.yml"
$ws.Cells.Item(3, 3).Value = 42
$ws.Cells.Item(3, 4).Value = 1000000
$ws.Cells.Item(3, 5).Value = 20
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 20
$ws.Cells.Item(3, 8).Value = "linear"
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = $false
$ws.Cells.Item(3, 11).Value = "none"
$ws.Cells.Item(3, 12).Value = 0

# ---------------------------------------------------------------------------
# Sheet 3: "Feature Details"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Feature Details")

# DS001's 20 features were actually discrete, not continuous.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 3).Value = "discrete"
}

# DS002 repeats the exact same per-feature stats as DS001, also discrete.
$features = @(
    @(0,  10,  2,   3.5),
    @(1,  -5,  1.5, -2.8),
    @(2,  0,   2.5,  0.1),
    @(3,  3,   1,    5.2),
    @(4,  -2,  0.8, -4.1),
    @(5,  15,  3,    1.7),
    @(6,  -8,  2.2, -0.05),
    @(7,  5,   1.8,  2.9),
    @(8,  -1,  0.9, -1.2),
    @(9,  12,  2.8,  0.8),
    @(10, -6,  1.3, -3.3),
    @(11, 8,   2.1,  0.4),
    @(12, -3,  1.1,  4.6),
    @(13, 20,  4,   -0.9),
    @(14, -10, 2.5,  2.1),
    @(15, 7,   1.6, -5),
    @(16, -4,  1.4,  0.3),
    @(17, 25,  5,    3.8),
    @(18, -12, 3.2, -1.6),
    @(19, 18,  3.8,  0)
)

$r = 22
foreach ($f in $features) {
    $ws.Cells.Item($r, 1).Value = "DS002"
    $ws.Cells.Item($r, 2).Value = "feature_" + $f[0]
    $ws.Cells.Item($r, 3).Value = "discrete"
    $ws.Cells.Item($r, 4).Value = $f[1]
    $ws.Cells.Item($r, 5).Value = $f[2]
    $ws.Cells.Item($r, 6).Value = $f[3]
    $ws.Cells.Item($r, 7).Value = $true
    $r++
}

# ---------------------------------------------------------------------------
# Sheet 4: "Weight Statistics"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Weight Statistics")

$ws.Cells.Item(3, 1).Value = "DS002"
$ws.Cells.Item(3, 2).Value = 20
$ws.Cells.Item(3, 3).Value = 11
$ws.Cells.Item(3, 4).Value = 8
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = -5
$ws.Cells.Item(3, 7).Value = 5.2
$ws.Cells.Item(3, 8).Value = 10.2

# ---------------------------------------------------------------------------
# Sheet 5: "File Metadata"
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("File Metadata")

$ws.Cells.Item(2, 2).Value = "configs\data_generation\dataset-config-001.yml"
$ws.Cells.Item(2, 3).Value = "data\n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
$ws.Cells.Item(2, 4).Value = "reports\figures\n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_plot.pdf"
$ws.Cells.Item(2, 5).Value = 106.1
$ws.Cells.Item(2, 6).Value = "f938398f95"
$ws.Cells.Item(2, 7).Value = "Linear function, 0C/20D features"

$ws.Cells.Item(3, 1).Value = "DS002"
$ws.Cells.Item(3, 2).Value = "configs\data_generation\dataset-config-002.yml"
$ws.Cells.Item(3, 3).Value = "data\n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_dataset.csv"
$ws.Cells.Item(3, 4).Value = "reports\figures\n1000000_f_init20_cont0_disc20_add0_pert-none_scl0_func-linear_noise0_plot.pdf"
$ws.Cells.Item(3, 5).Value = 106.1
$ws.Cells.Item(3, 6).Value = "f938398f95"
$ws.Cells.Item(3, 7).Value = "Linear function, 0C/20D features"

Write-Host "Edit complete"
